$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row and title-case municipality/state names
$ws.Range('A1').Value = 'mx_state'
$ws.Range('B1').Value = 'mx_municipality'
$ws.Range('C1').Value = 'n_matriculas'
$ws.Range('D1').Value = 'pct_matriculas'
$ws.Range('B5').Value = 'Pabellón De Arteaga'
$ws.Range('B6').Value = 'Rincón De Romos'
$ws.Range('B22').Value = 'Amatenango De La Frontera'
$ws.Range('B25').Value = 'Bejucal De Ocampo'
$ws.Range('B27').Value = 'Benemérito De Las Américas'
$ws.Range('B33').Value = 'Chiapa De Corzo'
$ws.Range('B37').Value = 'Comitán De Domínguez'
$ws.Range('B51').Value = 'Marqués De Comillas'
$ws.Range('B54').Value = 'Montecristo De Guerrero'
$ws.Range('B66').Value = 'Salto De Agua'
$ws.Range('B67').Value = 'San Cristóbal De Las Casas'
$ws.Range('B100').Value = 'Hidalgo Del Parral'
$ws.Range('B104').Value = 'San Francisco De Borja'
$ws.Range('A107').Value = 'Ciudad De México'
$ws.Range('A122').Value = 'Coahuila De Zaragoza'
$ws.Range('B144').Value = 'Pánuco De Coronado'
$ws.Range('A148').Value = 'Estado De México'
$ws.Range('B148').Value = 'Acambay De Ruíz Castañeda'
$ws.Range('B152').Value = 'Atizapán De Zaragoza'
$ws.Range('B156').Value = 'Chapa De Mota'
$ws.Range('B159').Value = 'Coacalco De Berriozábal'
$ws.Range('B164').Value = 'Ecatepec De Morelos'
$ws.Range('B167').Value = 'Ixtapan De La Sal'
$ws.Range('B174').Value = 'Naucalpan De Juárez'
$ws.Range('B179').Value = 'San Felipe Del Progreso'
$ws.Range('B180').Value = 'San Martín De Las Pirámides'
$ws.Range('B188').Value = 'Tlalnepantla De Baz'
$ws.Range('B192').Value = 'Valle De Bravo'
$ws.Range('B193').Value = 'Valle De Chalco Solidaridad'
$ws.Range('B196').Value = 'Villa Del Carbón'
$ws.Range('B201').Value = 'Apaseo El Alto'
$ws.Range('B202').Value = 'Apaseo El Grande'
$ws.Range('B209').Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Range('B217').Value = 'Purísima Del Rincón'
$ws.Range('B221').Value = 'San Diego De La Unión'
$ws.Range('B223').Value = 'San Francisco Del Rincón'
$ws.Range('B225').Value = 'San Luis De La Paz'
$ws.Range('B226').Value = 'San Miguel De Allende'
$ws.Range('B228').Value = 'Silao De La Victoria'
$ws.Range('B232').Value = 'Valle De Santiago'
$ws.Range('B236').Value = 'Acapulco De Juárez'
$ws.Range('B239').Value = 'Ajuchitlán Del Progreso'
$ws.Range('B240').Value = 'Alcozauca De Guerrero'
$ws.Range('B243').Value = 'Atoyac De Álvarez'
$ws.Range('B244').Value = 'Ayutla De Los Libres'
$ws.Range('B246').Value = 'Chilapa De Álvarez'
$ws.Range('B247').Value = 'Chilpancingo De Los Bravo'
$ws.Range('B250').Value = 'Coyuca De Benítez'
$ws.Range('B251').Value = 'Coyuca De Catalán'
$ws.Range('B256').Value = 'Huitzuco De Los Figueroa'
$ws.Range('B257').Value = 'Iguala De La Independencia'
$ws.Range('B269').Value = 'Tepecoacuilco De Trujano'
$ws.Range('B274').Value = 'Tlapa De Comonfort'
$ws.Range('B275').Value = 'Técpan De Galeana'
$ws.Range('B277').Value = 'Zihuatanejo De Azueta'
$ws.Range('B282').Value = 'Agua Blanca De Iturbide'
$ws.Range('B287').Value = 'Atotonilco El Grande'
$ws.Range('B292').Value = 'Cuautepec De Hinojosa'
$ws.Range('B296').Value = 'Huejutla De Reyes'
$ws.Range('B302').Value = 'Mixquiahuala De Juárez'
$ws.Range('B304').Value = 'Pachuca De Soto'
$ws.Range('B306').Value = 'Progreso De Obregón'
$ws.Range('B311').Value = 'Santiago De Anaya'
$ws.Range('B313').Value = 'Tenango De Doria'
$ws.Range('B315').Value = 'Tepehuacán De Guerrero'
$ws.Range('B316').Value = 'Tepeji Del Río De Ocampo'
$ws.Range('B317').Value = 'Tezontepec De Aldama'
$ws.Range('B320').Value = 'Tula De Allende'
$ws.Range('B321').Value = 'Tulancingo De Bravo'
$ws.Range('B323').Value = 'Zacualtipán De Ángeles'
$ws.Range('B327').Value = 'Autlán De Navarro'
$ws.Range('B336').Value = 'Lagos De Moreno'
$ws.Range('B342').Value = 'San Juanito De Escobedo'
$ws.Range('B345').Value = 'Tepatitlán De Morelos'
$ws.Range('B347').Value = 'Tlajomulco De Zúñiga'
$ws.Range('B351').Value = 'Unión De San Antonio'
$ws.Range('B352').Value = 'Yahualica De González Gallo'
$ws.Range('A355').Value = 'Michoacán De Ocampo'
$ws.Range('B401').Value = 'Puente De Ixtla'
$ws.Range('B402').Value = 'Tlaltizapán De Zapata'
$ws.Range('B415').Value = 'Mier Y Noriega'
$ws.Range('B418').Value = 'San Nicolás De Los Garza'
$ws.Range('B420').Value = 'Acatlán De Pérez Figueroa'
$ws.Range('B425').Value = 'Coicoyán De Las Flores'
$ws.Range('B428').Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Range('B429').Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Range('B430').Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Range('B431').Value = 'Huajuapan De León'
$ws.Range('B432').Value = 'Ixtlán De Juárez'
$ws.Range('B435').Value = 'Miahuatlán De Porfirio Díaz'
$ws.Range('B436').Value = 'Oaxaca De Juárez'
$ws.Range('B437').Value = 'Putla Villa De Guerrero'
$ws.Range('B442').Value = 'San Felipe Jalapa De Díaz'
$ws.Range('B446').Value = 'San José Del Progreso'
$ws.Range('B457').Value = 'San Mateo Del Mar'
$ws.Range('B460').Value = 'San Miguel Del Puerto'
$ws.Range('B490').Value = 'Tataltepec De Valdés'
$ws.Range('B491').Value = 'Teotitlán De Flores Magón'
$ws.Range('B492').Value = 'Villa De Tututepec'
$ws.Range('B493').Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Range('B500').Value = 'Chalchicomula De Sesma'
$ws.Range('B513').Value = 'Huitzilan De Serdán'
$ws.Range('B514').Value = 'Izúcar De Matamoros'
$ws.Range('B517').Value = 'Los Reyes De Juárez'
$ws.Range('B520').Value = 'Palmar De Bravo'
$ws.Range('B527').Value = 'San Salvador El Seco'
$ws.Range('B529').Value = 'Tecali De Herrera'
$ws.Range('B534').Value = 'Tetela De Ocampo'
$ws.Range('B549').Value = 'Amealco De Bonfil'
$ws.Range('B551').Value = 'Cadereyta De Montes'
$ws.Range('B556').Value = 'Jalpan De Serra'
$ws.Range('B558').Value = 'Pinal De Amoles'
$ws.Range('B561').Value = 'San Juan Del Río'
$ws.Range('B569').Value = 'Axtla De Terrazas'
$ws.Range('B572').Value = 'Cerro De San Pedro'
$ws.Range('B575').Value = 'Ciudad Del Maíz'
$ws.Range('B585').Value = 'San Ciro De Acosta'
$ws.Range('B590').Value = 'Santa María Del Río'
$ws.Range('B591').Value = 'Soledad De Graciano Sánchez'
$ws.Range('B597').Value = 'Villa De Arriaga'
$ws.Range('B598').Value = 'Villa De Ramos'
$ws.Range('B599').Value = 'Villa De Reyes'
$ws.Range('B634').Value = 'Soto La Marina'
$ws.Range('B640').Value = 'Apetatitlán De Antonio Carvajal'
$ws.Range('B643').Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Range('A648').Value = 'Veracruz De Ignacio De La Llave'
$ws.Range('B654').Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Range('B657').Value = 'Amatlán De Los Reyes'
$ws.Range('B665').Value = 'Cosamaloapan De Carpio'
$ws.Range('B666').Value = 'Cosautlán De Carvajal'
$ws.Range('B680').Value = 'Hueyapan De Ocampo'
$ws.Range('B681').Value = 'Ignacio De La Llave'
$ws.Range('B685').Value = 'Ixhuatlán De Madero'
$ws.Range('B686').Value = 'Ixhuatlán Del Café'
$ws.Range('B694').Value = 'Martínez De La Torre'
$ws.Range('B699').Value = 'Mixtla De Altamirano'
$ws.Range('B704').Value = 'Paso Del Macho'
$ws.Range('B706').Value = 'Poza Rica De Hidalgo'
$ws.Range('B714').Value = 'Sayula De Alemán'
$ws.Range('B715').Value = 'Soledad De Doblado'
$ws.Range('B717').Value = 'Tatahuicapan De Juárez'
$ws.Range('B737').Value = 'Vega De Alatorre'
$ws.Range('B744').Value = 'Zozocolco De Hidalgo'
$ws.Range('B762').Value = 'Nochistlán De Mejía'
$ws.Range('A775').Value = 'Total'

# Remove trailing metadata rows (777:781); dimension auto-shrinks to A1:D775
$ws.Rows("777:781").Delete()
